$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2 (123 Rue Principale...) into the new row 12, then change the address
$ws.Range("A2:I2").Copy()
$ws.Range("A12:I12").PasteSpecial()
$excel.CutCopyMode = $false
$ws.Range("A12").Value = "725 Rue Jolliet, Saguenay, QC, G7J 2P7"

# Widen column A so the full addresses are visible
$ws.Range("A1").EntireColumn.ColumnWidth = 55

# Move the active selection to the next empty row, as left by the editor
$ws.Range("A13").Select()
